# Quarterly income statement update:
#  - drop the stale "6 ماهه منتهی به 1399/06" period column (old column D)
#  - append the new "12 ماهه منتهی به 1401/12" period column (new column M)
#  - refresh the "12 ماهه منتهی به 1400/12" column (now column I) with
#    recalculated values (read_price algorithm change)
#  - refresh the "تاریخ انتشار" (publish date) caption for that column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "wide" column width (the 29-char columns used for every 3rd
# period) before the shift, so the newly appended column M can match it.
$wideColumnWidth = $ws.Range("F1").ColumnWidth

# Remove the obsolete first period column; every remaining column (and its
# data) shifts one slot to the left automatically.
$ws.Columns("D").Delete()

# New column M is the 3rd column in its group of 3 (K,L,M), so it gets the
# wide width just like D/G/J used to.
$ws.Range("M1").ColumnWidth = $wideColumnWidth

# --- New period header / publish date for the freshly shifted-in column I ---
$ws.Range("I9").Value = "1402-02-12 (10)"

# --- Recalculated values for column I ("12 ماهه منتهی به 1400/12") ---
$ws.Range("I12").Value = -4982668
$ws.Range("I13").Value = 3183455
$ws.Range("I14").Value = -439067
$ws.Range("I16").Value = -344346
$ws.Range("I17").Value = 2400042
$ws.Range("I20").Value = 2637540
$ws.Range("I22").Value = 2296297
$ws.Range("I24").Value = 2296297
$ws.Range("I25").Value = 2088
$ws.Range("I27").Value = 2088

# --- New column M: "12 ماهه منتهی به 1401/12" period ---
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-12 (2)"

$ws.Range("M11").Value = 12146108
$ws.Range("M12").Value = -7538475
$ws.Range("M13").Value = 4607633
$ws.Range("M14").Value = -744022
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 113993
$ws.Range("M17").Value = 3977604
$ws.Range("M18").Value = -7373
$ws.Range("M19").Value = 492256
$ws.Range("M20").Value = 4462487
$ws.Range("M21").Value = -409181
$ws.Range("M22").Value = 4053306
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 4053306
$ws.Range("M25").Value = 3685
$ws.Range("M26").Value = 1100000
$ws.Range("M27").Value = 3685
